$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at row 11, pushing old rows 11-12 (Buy/Sell Decision, Tab 2)
# down to rows 17-18. Using EntireRow.Insert() (rather than writing over existing
# cells) preserves the row-11 wrapped-text height (ht="165") on the moved row
# without forcing a "customHeight" flag, matching native Excel row-insert behavior.
$ws.Range("A11:A16").EntireRow.Insert()

# The moved rows used to be tagged "Tab 2"; they are now "Tab 3".
$ws.Range("A17").Value = "Tab 3"
$ws.Range("A18").Value = "Tab 3"

# Populate the newly inserted rows 11-16 with the new "Player Stats" /
# "Player Rating" configuration rows (Tab 2).
$ws.Range("A11").Value = "Tab 2"
$ws.Range("B11").Value = "Player Stats"
$ws.Range("C11").Value = "No. of clubs played for"
$ws.Range("D11").Value = "Different clubs in Europe's top 5 leagues"
$ws.Range("E11").Value = "y"
$ws.Range("F11").Value = "numeric"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "Comment on this"
$ws.Range("K11").Value = "y"
$ws.Range("L11").Value = "text"
$ws.Range("N11").Value = "Comments"
$ws.Range("P11").Value = "n"

$ws.Range("A12").Value = "Tab 2"
$ws.Range("B12").Value = "Player Stats"
$ws.Range("C12").Value = "Goals scored for country"
$ws.Range("D12").Value = "Goals scored by the player for the country"
$ws.Range("E12").Value = "y"
$ws.Range("F12").Value = "numeric"
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = "n"
$ws.Range("P12").Value = "n"

$ws.Range("A13").Value = "Tab 2"
$ws.Range("B13").Value = "Player Stats"
$ws.Range("C13").Value = "Goals scored for club"
$ws.Range("D13").Value = "Goals scored by the player for all their clubs"
$ws.Range("E13").Value = "y"
$ws.Range("F13").Value = "numeric"
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = "n"
$ws.Range("P13").Value = "n"

$ws.Range("A14").Value = "Tab 2"
$ws.Range("B14").Value = "Player Stats"
$ws.Range("C14").Value = "Goals scored in youth career"
$ws.Range("D14").Value = "Goals scored by the player in their youth career"
$ws.Range("E14").Value = "y"
$ws.Range("F14").Value = "numeric"
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = "n"
$ws.Range("P14").Value = "n"

$ws.Range("A15").Value = "Tab 2"
$ws.Range("B15").Value = "Player Rating"
$ws.Range("C15").Value = "Player country score"
$ws.Range("E15").Value = "e"
$ws.Range("F15").Value = "read-only"
$ws.Range("H15").Value = "Goals scored for country"

$ws.Range("A16").Value = "Tab 2"
$ws.Range("B16").Value = "Player Rating"
$ws.Range("C16").Value = "Player club score"
$ws.Range("E16").Value = "e"
$ws.Range("F16").Value = "read-only"
$ws.Range("H16").Value = "Goals scored for club"

# Resize the table and autofilter to cover the new data range.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:S18"))

# Move the selection to match the edited workbook's saved cursor position.
$ws.Range("H16").Select() | Out-Null
